$wb = $excel.ActiveWorkbook

$changes = @{
    2  = 11628
    4  = 7
    5  = 1059
    7  = 76
    8  = 49
    10 = 11003
    11 = 4236
    14 = 13
    17 = 69
    20 = 468
    21 = 11182
    22 = 11003
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $changes.Keys) {
        $ws.Range("F$row").Value = $changes[$row]
    }
}
